$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 16206.257
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 16206.257
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 48618.771
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -48954.771
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 500
$ws.Range("K18").Value = 500
$ws.Range("M18").Value = -216
$ws.Range("H28").Value = 22033.945
$ws.Range("I28").Value = 30319.584
$ws.Range("J28").Value = 5462.6665
$ws.Range("K28").Value = 30319.584
$ws.Range("L28").Value = 5462.6665
$ws.Range("M28").Value = -29834.584
$ws.Range("N28").Value = -6432.6665
$ws.Range("H43").Value = 2998.5
$ws.Range("I43").Value = 2495
$ws.Range("K43").Value = 2495
$ws.Range("M43").Value = -2426
$ws.Range("H62").Value = 1239.625
$ws.Range("I62").Value = 1014.44446
$ws.Range("J62").Value = 1529.1428
$ws.Range("K62").Value = 1014.44446
$ws.Range("L62").Value = 1529.1428
$ws.Range("M62").Value = -390.44446
$ws.Range("N62").Value = -2777.1428
$ws.Range("H65").Value = 1239.625
$ws.Range("I65").Value = 1014.44446
$ws.Range("J65").Value = 1529.1428
$ws.Range("K65").Value = 5072.2223
$ws.Range("L65").Value = 7645.714
$ws.Range("M65").Value = -1952.2223
$ws.Range("N65").Value = -13885.714
$ws.Range("H132").Value = 4019.6667
$ws.Range("I132").Value = 4097.6665
$ws.Range("J132").Value = 3883.1667
$ws.Range("K132").Value = 12292.9995
$ws.Range("L132").Value = 11649.5001
$ws.Range("M132").Value = -9762.999500000002
$ws.Range("N132").Value = -16709.5001
$ws.Range("H137").Value = 3411.6736
$ws.Range("I137").Value = 859.3929000000001
$ws.Range("J137").Value = 6814.7144
$ws.Range("K137").Value = 2578.1787
$ws.Range("L137").Value = 20444.1432
$ws.Range("M137").Value = -28.17870000000039
$ws.Range("N137").Value = -25544.1432
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 11170.625
$ws.Range("I45").Value = 14050.167
$ws.Range("J45").Value = 2532
$ws.Range("K45").Value = 14050.167
$ws.Range("L45").Value = 2532
$ws.Range("M45").Value = -13673.167
$ws.Range("N45").Value = -3286
$ws.Range("H61").Value = 3790.75
$ws.Range("I61").Value = 3711.2727
$ws.Range("J61").Value = 3965.6
$ws.Range("K61").Value = 3711.2727
$ws.Range("L61").Value = 3965.6
$ws.Range("M61").Value = -3499.2727
$ws.Range("N61").Value = -4389.6
$ws.Range("H74").Value = 5738.727
$ws.Range("I74").Value = 973.06665
$ws.Range("J74").Value = 15950.857
$ws.Range("K74").Value = 973.06665
$ws.Range("L74").Value = 15950.857
$ws.Range("M74").Value = -99.06664999999998
$ws.Range("N74").Value = -17698.857
$ws.Range("H77").Value = 5738.727
$ws.Range("I77").Value = 973.06665
$ws.Range("J77").Value = 15950.857
$ws.Range("K77").Value = 4865.33325
$ws.Range("L77").Value = 79754.285
$ws.Range("M77").Value = -497.3332499999997
$ws.Range("N77").Value = -88490.285
$ws.Range("H88").Value = 2450
$ws.Range("I88").Value = 1800
$ws.Range("J88").Value = 2666.6667
$ws.Range("K88").Value = 1800
$ws.Range("L88").Value = 2666.6667
$ws.Range("M88").Value = -1394
$ws.Range("N88").Value = -3478.6667
$ws.Range("H91").Value = 2450
$ws.Range("I91").Value = 1800
$ws.Range("J91").Value = 2666.6667
$ws.Range("K91").Value = 1800
$ws.Range("L91").Value = 2666.6667
$ws.Range("M91").Value = -396
$ws.Range("N91").Value = -5474.6667
$ws.Range("H132").Value = 5270.5425
$ws.Range("I132").Value = 5747.875
$ws.Range("J132").Value = 4943.2285
$ws.Range("K132").Value = 17243.625
$ws.Range("L132").Value = 14829.6855
$ws.Range("M132").Value = -14713.625
$ws.Range("N132").Value = -19889.6855
$ws.Range("H136").Value = 3790.75
$ws.Range("I136").Value = 3711.2727
$ws.Range("J136").Value = 3965.6
$ws.Range("K136").Value = 11133.8181
$ws.Range("L136").Value = 11896.8
$ws.Range("M136").Value = -8583.8181
$ws.Range("N136").Value = -16996.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2412.8
$ws.Range("I86").Value = 1926.6666
$ws.Range("J86").Value = 2686.25
$ws.Range("K86").Value = 1926.6666
$ws.Range("L86").Value = 2686.25
$ws.Range("M86").Value = -803.6666
$ws.Range("N86").Value = -4932.25
$ws.Range("H89").Value = 2412.8
$ws.Range("I89").Value = 1926.6666
$ws.Range("J89").Value = 2686.25
$ws.Range("K89").Value = 9633.333000000001
$ws.Range("L89").Value = 13431.25
$ws.Range("M89").Value = -4017.333000000001
$ws.Range("N89").Value = -24663.25
$ws.Range("H99").Value = 1466.6923
$ws.Range("I99").Value = 980.625
$ws.Range("J99").Value = 2244.4
$ws.Range("K99").Value = 980.625
$ws.Range("L99").Value = 2244.4
$ws.Range("M99").Value = 517.375
$ws.Range("N99").Value = -5240.4
$ws.Range("H134").Value = 1571.7407
$ws.Range("I134").Value = 1015.9524
$ws.Range("J134").Value = 3517
$ws.Range("K134").Value = 3047.8572
$ws.Range("L134").Value = 10551
$ws.Range("M134").Value = -512.8571999999999
$ws.Range("N134").Value = -15621

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17878976
$ws.Range("I31").Value = 41668176
$ws.Range("K31").Value = 41668176
$ws.Range("M31").Value = -41667881
$ws.Range("H34").Value = 17878976
$ws.Range("I34").Value = 41668176
$ws.Range("K34").Value = 41668176
$ws.Range("M34").Value = -41667974
$ws.Range("H58").Value = 1287.6552
$ws.Range("I58").Value = 1210.3636
$ws.Range("K58").Value = 1210.3636
$ws.Range("M58").Value = -1007.3636
$ws.Range("H132").Value = 35720920
$ws.Range("I132").Value = 90921620
$ws.Range("J132").Value = 2820.5293
$ws.Range("K132").Value = 272764860
$ws.Range("L132").Value = 8461.5879
$ws.Range("M132").Value = -272762330
$ws.Range("N132").Value = -13521.5879
$ws.Range("H134").Value = 1130.1538
$ws.Range("I134").Value = 898.125
$ws.Range("K134").Value = 2694.375
$ws.Range("M134").Value = -159.375
$ws.Range("H136").Value = 1287.6552
$ws.Range("I136").Value = 1210.3636
$ws.Range("K136").Value = 3631.0908
$ws.Range("M136").Value = -1081.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 795.34784
$ws.Range("I7").Value = 62
$ws.Range("J7").Value = 999.05554
$ws.Range("K7").Value = 186
$ws.Range("L7").Value = 2997.16662
$ws.Range("M7").Value = -74
$ws.Range("N7").Value = -3221.16662
$ws.Range("H137").Value = 10211317
$ws.Range("I137").Value = 184831.67
$ws.Range("K137").Value = 554495.01
$ws.Range("M137").Value = -549395.01

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2905.6365
$ws.Range("I102").Value = 3314.6667
$ws.Range("K102").Value = 3314.6667
$ws.Range("M102").Value = -1692.6667
$ws.Range("H122").Value = 957
$ws.Range("I122").Value = 914
$ws.Range("K122").Value = 2742
$ws.Range("M122").Value = -292
$ws.Range("H132").Value = 7768.0835
$ws.Range("I132").Value = 8733.474
$ws.Range("J132").Value = 4099.6
$ws.Range("K132").Value = 26200.422
$ws.Range("L132").Value = 12298.8
$ws.Range("M132").Value = -23670.422
$ws.Range("N132").Value = -17358.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7150.08
$ws.Range("I132").Value = 8547.444
$ws.Range("J132").Value = 3556.8572
$ws.Range("K132").Value = 25642.332
$ws.Range("L132").Value = 10670.5716
$ws.Range("M132").Value = -23112.332
$ws.Range("N132").Value = -15730.5716
$ws.Range("H136").Value = 6081.423
$ws.Range("I136").Value = 2400.6667
$ws.Range("J136").Value = 14363.125
$ws.Range("K136").Value = 7202.000100000001
$ws.Range("L136").Value = 43089.375
$ws.Range("M136").Value = -4652.000100000001
$ws.Range("N136").Value = -48189.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 11112618
$ws.Range("I122").Value = 15385774
$ws.Range("J122").Value = 2412
$ws.Range("K122").Value = 46157322
$ws.Range("L122").Value = 7236
$ws.Range("M122").Value = -46154872
$ws.Range("N122").Value = -12136
$ws.Range("H132").Value = 25007378
$ws.Range("I132").Value = 41677256
$ws.Range("J132").Value = 2564.125
$ws.Range("K132").Value = 125031768
$ws.Range("L132").Value = 7692.375
$ws.Range("M132").Value = -125029238
$ws.Range("N132").Value = -12752.375
$ws.Range("H136").Value = 3849.3948
$ws.Range("I136").Value = 8910.083000000001
$ws.Range("J136").Value = 1513.6923
$ws.Range("K136").Value = 26730.249
$ws.Range("L136").Value = 4541.0769
$ws.Range("M136").Value = -24180.249
$ws.Range("N136").Value = -9641.0769
